# dmcar_serialisation_template.xlsx edit
# Commit: "minor refactor loader content to graphloader.py, ontology tweaks"
#
# Content change: the "Schema" sheet is a self-describing metadata table
# (it lists the attributes of the "Schema" schema itself, in the
# SchemaSchema block). The "ParentNamespace" attribute row (row 5) is
# removed from that table - an ontology tweak dropping a now-unused
# schema attribute. Deleting the row shifts all subsequent rows up by
# one (their Sequence numbers are literal values and are NOT
# renumbered, which is exactly what a native row delete does).

$wb = $excel.ActiveWorkbook

$wsSchema = $wb.Worksheets.Item("Schema")
$wsSchema.Rows.Item(5).Delete()

# View-state touch-ups so the saved selection / active sheet matches the
# post-edit workbook: Mapping's cursor moved to I2, Schema's cursor moved
# to H17 (Schema is, and remains, the active tab), DMCAR and Requirement
# keep their prior cursors (S1 / B8 respectively).
$wsMapping = $wb.Worksheets.Item("Mapping")
$wsMapping.Activate()
$wsMapping.Range("I2").Select()

$wsDMCAR = $wb.Worksheets.Item("DMCAR")
$wsDMCAR.Activate()
$wsDMCAR.Range("S1").Select()

$wsRequirement = $wb.Worksheets.Item("Requirement")
$wsRequirement.Activate()
$wsRequirement.Range("B8").Select()

$wsSchema.Activate()
$wsSchema.Range("H17").Select()
